# Add PGC* to 3 studies
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Ripke (PGC1)"
$ws.Range("A5").Value = "Wray (PGC2)"
$ws.Range("A8").Value = "Adams (PGC3)"

$ws.Range("A9").Select()
